# Adicionando graficos de estatisticas descritivas e analises de profile
# This updates the dataset sheet: the old "m5" metric column (E) is dropped
# (shifting m6/repository/version/.../Reliability/total one column to the
# left) and every row is refreshed with newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column E ("m5"); Excel shifts F:N left to E:M and the dimension
# becomes A1:M8 automatically. The header row ends up correct as a result
# (E1 becomes "m6", F1 becomes "repository", ..., M1 becomes "total").
$ws.Columns("E").Delete()

# Refreshed data rows (columns A:M) for rows 2-8.
$data = @(
    @(1, 0.3333333333333333, 1, 1, 1, "Archives", "03-10-2021-15-47", 239, 0.77, 1, 0.385, 0.5, 0.885),
    @(1, 0.3333333333333333, 1, 1, 0.8333333333333334, "Archives", "13-09-2021-15-00", 109, 0.77, 0.8833333333333333, 0.385, 0.4416666666666667, 0.8266666666666667),
    @(1, 0.3333333333333333, 1, 1, 0.8333333333333334, "Profile", "03-10-2021-15-48", 112, 0.77, 0.8833333333333333, 0.385, 0.4416666666666667, 0.8266666666666667),
    @(1, 0.3333333333333333, 1, 1, 0.8333333333333334, "Profile", "13-09-2021-14-00", 112, 0.77, 0.8833333333333333, 0.385, 0.4416666666666667, 0.8266666666666667),
    @(0.88, 0.04, 1, 1, 0.92, "Frontend", "03-10-2021-15-49", 1417, 0.6335999999999999, 0.944, 0.3168, 0.472, 0.7887999999999999),
    @(0.5, 0.25, 1, 1, 0.5, "Frontend", "13-09-2021-20-00", 48, 0.5775, 0.6499999999999999, 0.28875, 0.325, 0.61375),
    @(0.9423076923076923, 0.01923076923076923, 0.9230769230769231, 1, 0.9038461538461539, "Frontend", "17-10-2021-15-30", 4874, 0.6219230769230769, 0.9326923076923077, 0.3109615384615385, 0.4663461538461539, 0.7773076923076923)
)

$columns = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $columns.Length; $j++) {
        $cellRef = $columns[$j] + $rowNum
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}
